$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.104873180389404
$ws.Range("B1").Value = 1.717735767364502
$ws.Range("C1").Value = 9.148412704467773
$ws.Range("D1").Value = 2.39408802986145
$ws.Range("E1").Value = 1.2852543592453
